$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws2 = $wb.Worksheets.Item("演出")
$ws3 = $wb.Worksheets.Item("本地生活")
$ws4 = $wb.Worksheets.Item("全部类型")

$ws1.Range("F3").Value = 897
$ws1.Range("F5").Value = 2304
$ws1.Range("F6").Value = 1372
$ws1.Range("F8").Value = 818
$ws1.Range("F9").Value = 1185
$ws1.Range("F11").Value = 3077
$ws1.Range("F12").Value = 36
$ws1.Range("F14").Value = 1108
$ws1.Range("F19").Value = 1157
$ws1.Range("F20").Value = 1157
$ws1.Range("F25").Value = 267
$ws1.Range("F29").Value = 847
$ws1.Range("F32").Value = 80
$ws1.Range("F33").Value = 1065
$ws1.Range("F34").Value = 5083
$ws1.Range("F36").Value = 285
$ws2.Range("F6").Value = 416
$ws2.Range("F22").Value = 305
$ws2.Range("F24").Value = 53
$ws2.Range("F25").Value = 395
$ws2.Range("F28").Value = 699
$ws3.Range("F6").Value = 433
$ws4.Range("F6").Value = 897
$ws4.Range("F7").Value = 416
$ws4.Range("F9").Value = 2304
$ws4.Range("F10").Value = 1372
$ws4.Range("F12").Value = 818
$ws4.Range("F13").Value = 1185
$ws4.Range("F16").Value = 3077
$ws4.Range("F17").Value = 36
$ws4.Range("F20").Value = 1108
$ws4.Range("F21").Value = 433
$ws4.Range("C24").Value = "上海·林肯爵士乐上海中心大乐队·经典动漫二次元ACG爵士音乐会"
$ws4.Range("D24").Value = "南京东路139号4层（近江西中路） 林肯爵士乐上海中心"
$ws4.Range("E24").Value = "2024.04.27 20:00-04.27 21:30"
$ws4.Range("F24").Value = 4
$ws4.Range("G24").Value = 220
$ws4.Range("H24").Value = "https://show.bilibili.com/platform/detail.html?id=83976"
$ws4.Range("I24").Value = "//i1.hdslb.com/bfs/openplatform/202404/92F3IJmA1712654498838.jpeg"
$ws4.Range("B25").Value = "'2024-04-27"
$ws4.Range("C25").Value = "上海·第五十六届燃梦星辰动漫嘉年华"
$ws4.Range("D25").Value = "年家浜东路238号 上海绿地缤纷广场(周浦店)"
$ws4.Range("E25").Value = "2024.04.27 10:30-04.27 16:30"
$ws4.Range("F25").Value = 616
$ws4.Range("G25").Value = 58.8
$ws4.Range("H25").Value = "https://show.bilibili.com/platform/detail.html?id=83814"
$ws4.Range("I25").Value = "//i2.hdslb.com/bfs/openplatform/202404/um6MUtv61712460652109.jpeg"
$ws4.Range("F26").Value = 1157
$ws4.Range("C27").Value = "上海·HATSUNE MIKU meets niko and ... 集章之旅"
$ws4.Range("D27").Value = "淮海中路775号 niko and ......"
$ws4.Range("E27").Value = "2024.05.01 10:00-06.02 22:00"
$ws4.Range("F27").Value = 1157
$ws4.Range("G27").Value = 46
$ws4.Range("H27").Value = "https://show.bilibili.com/platform/detail.html?id=83163"
$ws4.Range("I27").Value = "//i2.hdslb.com/bfs/openplatform/202403/9lMpza7M1711528161190.jpeg"
$ws4.Range("C28").Value = "上海·coser动漫展03"
$ws4.Range("D28").Value = "海潮路133号B1 JUMP工坊"
$ws4.Range("E28").Value = "2024.05.01 10:00-05.02 17:00"
$ws4.Range("F28").Value = 176
$ws4.Range("G28").Value = 60
$ws4.Range("H28").Value = "https://show.bilibili.com/platform/detail.html?id=83563"
$ws4.Range("I28").Value = "//i0.hdslb.com/bfs/openplatform/202403/1avjQYmV1711705372433.jpeg"
$ws4.Range("C29").Value = "上海·「Azurock2.0」Azurose ACG Cover Live跨次元乐队音乐现场"
$ws4.Range("D29").Value = "愚园路1398号B1层B1-02 育音堂(音乐公园店)"
$ws4.Range("E29").Value = "2024.05.01 18:30-05.01 21:00"
$ws4.Range("F29").Value = 46
$ws4.Range("G29").Value = 88
$ws4.Range("H29").Value = "https://show.bilibili.com/platform/detail.html?id=83763"
$ws4.Range("I29").Value = "//i0.hdslb.com/bfs/openplatform/202404/hHTjtr5h1711980928508.jpeg"
$ws4.Range("C30").Value = "上海·第五十七届燃梦星辰动漫嘉年华"
$ws4.Range("D30").Value = "云锦路500号(近11号线地铁站5号口) 绿地滨江CLUB"
$ws4.Range("E30").Value = "2024.05.01 10:30-05.01 16:30"
$ws4.Range("F30").Value = 548
$ws4.Range("G30").Value = 58.8
$ws4.Range("H30").Value = "https://show.bilibili.com/platform/detail.html?id=83807"
$ws4.Range("I30").Value = "//i2.hdslb.com/bfs/openplatform/202404/RGLpPX211712156496032.jpeg"
$ws4.Range("F32").Value = 267
$ws4.Range("F33").Value = 53
$ws4.Range("F36").Value = 395
$ws4.Range("F37").Value = 699
$ws4.Range("F38").Value = 847
$ws4.Range("C41").Value = "上海·集训试炼·排球少年only"
$ws4.Range("D41").Value = "顾村镇蕰川路6号 智慧湾科创园"
$ws4.Range("E41").Value = "2024.05.05 10:00-05.05 18:00"
$ws4.Range("F41").Value = 80
$ws4.Range("G41").Value = 79
$ws4.Range("H41").Value = "https://show.bilibili.com/platform/detail.html?id=83986"
$ws4.Range("I41").Value = "//i0.hdslb.com/bfs/openplatform/202404/v13YTZum1712580327011.jpeg"
$ws4.Range("B42").Value = "'2024-05-18"
$ws4.Range("C42").Value = "上海·S·CGE动漫游戏嘉年华"
$ws4.Range("D42").Value = "军工路1076号 纪希片场(秀场)"
$ws4.Range("E42").Value = "2024.05.18 10:00-05.19 17:00"
$ws4.Range("F42").Value = 1065
$ws4.Range("G42").Value = 70
$ws4.Range("H42").Value = "https://show.bilibili.com/platform/detail.html?id=81204"
$ws4.Range("I42").Value = "//i0.hdslb.com/bfs/openplatform/202403/B4thpI0H1711678623426.jpeg"
$ws4.Range("C43").Value = "上海·原神×崩坏×星铁only旅行盛宴2.0"
$ws4.Range("D43").Value = "西藏南路1号 上海大世界"
$ws4.Range("F43").Value = 5083
$ws4.Range("G43").Value = 65
$ws4.Range("H43").Value = "https://show.bilibili.com/platform/detail.html?id=81276"
$ws4.Range("I43").Value = "//i0.hdslb.com/bfs/openplatform/202403/uqKEdQrT1710486553826.jpeg"
$ws4.Range("C44").Value = "上海·松本大辉2024上海粉丝见面会"
$ws4.Range("D44").Value = "上海市普陀区宜昌路179号 万代南梦宫上海文化中心"
$ws4.Range("E44").Value = "2024.05.18 14:00-05.18 20:30"
$ws4.Range("F44").Value = 63
$ws4.Range("G44").Value = 380
$ws4.Range("H44").Value = "https://show.bilibili.com/platform/detail.html?id=83562"
$ws4.Range("I44").Value = "//i2.hdslb.com/bfs/openplatform/202404/sTcvzhSb1712052549414.jpeg"
$ws4.Range("C45").Value = "上海·第五人格ONLY2.0"
$ws4.Range("D45").Value = "顾村镇蕰川路6号 智慧湾科创园"
$ws4.Range("E45").Value = "2024.05.18 10:00-05.18 17:00"
$ws4.Range("F45").Value = 543
$ws4.Range("G45").Value = 60
$ws4.Range("H45").Value = "https://show.bilibili.com/platform/detail.html?id=83269"
$ws4.Range("I45").Value = "//i1.hdslb.com/bfs/openplatform/202403/dG2qpfiw1711167409798.jpeg"
$ws4.Range("B46").Value = "'2024-05-19"
$ws4.Range("C46").Value = "上海·优声之形系列活动之内田彩见面会"
$ws4.Range("D46").Value = "宜昌路179号 万代南梦宫上海文化中心"
$ws4.Range("E46").Value = "2024.05.19 14:00-05.19 15:30"
$ws4.Range("F46").Value = 453
$ws4.Range("G46").Value = 380
$ws4.Range("H46").Value = "https://show.bilibili.com/platform/detail.html?id=83741"
$ws4.Range("I46").Value = "//i2.hdslb.com/bfs/openplatform/202404/6GYq69Sq1712128553852.jpeg"
$ws4.Range("C47").Value = "上海·恋与深空×恋与制作人only"
$ws4.Range("D47").Value = "顾村镇蕰川路6号 智慧湾科创园"
$ws4.Range("E47").Value = "2024.05.19 10:00-05.19 17:00"
$ws4.Range("F47").Value = 285
$ws4.Range("G47").Value = 60
$ws4.Range("H47").Value = "https://show.bilibili.com/platform/detail.html?id=83598"
$ws4.Range("I47").Value = "//i2.hdslb.com/bfs/openplatform/202403/jskKqUvJ1711165688442.jpeg"
$ws4.Range("B48").Value = "'2024-06-01"
$ws4.Range("C48").Value = "上海·英雄时代2024哈瓦西钢琴演奏会"
$ws4.Range("D48").Value = "高青西路777号 上海前滩31演艺中心"
$ws4.Range("E48").Value = "2024.06.01 19:30-06.01 21:00"
$ws4.Range("F48").Value = 1
$ws4.Range("G48").Value = 499
$ws4.Range("H48").Value = "https://show.bilibili.com/platform/detail.html?id=83893"
$ws4.Range("I48").Value = "//i1.hdslb.com/bfs/openplatform/202404/pHVcjZyP1712566658767.jpeg"
$ws4.Range("B49").Value = "'2024-06-22"
$ws4.Range("C49").Value = "上海·「多厨狂喜」白金交响乐团二次元交响音乐会"
$ws4.Range("D49").Value = "丁香路425号 上海东方艺术中心"
$ws4.Range("E49").Value = "2024.06.22 19:30-06.22 21:30"
$ws4.Range("F49").Value = 760
$ws4.Range("G49").Value = 188
$ws4.Range("H49").Value = "https://show.bilibili.com/platform/detail.html?id=82731"
$ws4.Range("I49").Value = "//i0.hdslb.com/bfs/openplatform/202403/K3AlF8sr1710230449280.jpeg"
$ws4.Range("B50").Value = "'2024-07-12"
$ws4.Range("C50").Value = "上海·魔都COS Voice暑假漫展"
$ws4.Range("D50").Value = "长寿路街道澳门路168号 月星家居"
$ws4.Range("E50").Value = "2024.07.12 10:00-07.14 16:00"
$ws4.Range("F50").Value = 10
$ws4.Range("G50").Value = 49
$ws4.Range("H50").Value = "https://show.bilibili.com/platform/detail.html?id=84008"
$ws4.Range("I50").Value = "//i2.hdslb.com/bfs/openplatform/202404/jldYwFXi1712678080397.jpeg"
$ws4.Range("B51").Value = "'2024-07-19"
$ws4.Range("C51").Value = "上海·《你的名字》《天气之子》《铃芽之旅》——新海诚动漫三部曲钢琴演奏会"
$ws4.Range("D51").Value = "丁香路425号(上海科技馆地铁站1号口步行460米) 上海东方艺术中心音乐厅"
$ws4.Range("E51").Value = "2024.07.19 19:30-07.19 21:30"
$ws4.Range("F51").Value = 44
$ws4.Range("G51").Value = 80
$ws4.Range("H51").Value = "https://show.bilibili.com/platform/detail.html?id=83479"
$ws4.Range("I51").Value = "//i1.hdslb.com/bfs/openplatform/202403/GpyueuYA1711508106584.jpeg"
